# Loan RBI, Variable Instalments
#
# Inserts a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet - shifting the old N/O/P columns one slot to
# the right (N->O, O->P, P->Q) - and moves the active selection/tab back
# onto the "Repayment schedule" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column at position 14 (N), pushing existing N/O/P -> O/P/Q.
$ws.Columns.Item(14).Insert() | Out-Null

# Give the freshly inserted column the same (manually set) width the
# author gave it - matches column M's width.
$ws.Columns.Item(14).ColumnWidth = 9.83

# Move the selection/active sheet back to the repayment schedule sheet
# (this also flips the workbook's tracked active tab and the
# tabSelected flag away from the sheet that used to be active).
$ws.Activate() | Out-Null
$ws.Range("R6").Select() | Out-Null
